# Updated capital structure database
# Refreshes the computed financial ratios for the Pakistan / Air Transport
# rows (2 and 3) of the "earnings_debt" sheet to reflect the latest
# capital-structure inputs. Columns W (roe) and Y (roe_cost_equity) are
# newly populated; all other touched columns are recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "D"  = 0.0278
    "G"  = 0.003998919211024047
    "H"  = 0.003998919211024047
    "I"  = -0.1128073493650365
    "J"  = -0.1128073493650365
    "K"  = -323.7
    "L"  = -0.4373142393947581
    "U"  = 67.59999999999999
    "V"  = 0.4369747899159664
    "W"  = 0.1689986425811841
    "X"  = 1.019985182130581
    "Y"  = -0.850986539549397
    "Z"  = 15.42083333333336
    "AA" = -1.739583333333336
    "AB" = 0.119235985636019
    "AC" = -1.858819318969355
    "AD" = 2140.4
    "AE" = 0
    "AF" = 2140.4
    "AG" = 2072.8
    "AH" = 0.9325955296065531
    "AI" = -22.48319327731097
    "AJ" = 0.9305499438832773
    "AK" = -12.73218673218675
    "AL" = 210.9
    "AM" = 210.291
    "AN" = -68.38338658146965
    "AO" = -0.3959222380275011
    "AP" = -66.22364217252397
    "AQ" = -0.3970688236776657
}

foreach ($row in 2..3) {
    foreach ($col in $newValues.Keys) {
        $ws.Range("$col$row").Value2 = $newValues[$col]
    }
}
